# Minor corrections to a couple of code-sample slides.

$p = $ppt.ActivePresentation

# --- Slide 74: "ErrorHandler.getInstance().reportError(e);" -> "errorHandler.reportError(e);" ---
$s74 = $p.Slides.Item(74)
$shp74 = $s74.Shapes.Item(2)
$tr74 = $shp74.TextFrame.TextRange
$para74 = $tr74.Paragraphs(10, 1)

# Remove the now-unneeded "()." + "reportError" runs (kept text merges into the
# first run), then rewrite the first run's text as "errorHandler.reportError".
$mid74 = $para74.Characters(33, 14)
$mid74.Text = ""
$first74 = $para74.Characters(9, 24)
$first74.Text = "errorHandler.reportError"

# --- Slide 76: merge the trailing whitespace run into the "StandardCharsets.UTF_8));" run ---
$s76 = $p.Slides.Item(76)
$shp76 = $s76.Shapes.Item(2)
$tr76 = $shp76.TextFrame.TextRange
$para76 = $tr76.Paragraphs(2, 1)

# Delete the separate whitespace-only run (its formatting is dropped), then
# prepend the same whitespace onto the following run so it keeps that run's
# own formatting (dirty="0").
$ws76 = $para76.Characters(61, 37)
$ws76.Text = ""
$tail76 = $para76.Characters(61, 25)
$tail76.Text = "                                     StandardCharsets.UTF_8));"
